$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log entry (row 6): Alumno, Tarea, Fecha
$ws.Range("B6").Value = "Gaizka"
$ws.Range("C6").Value = "Update hecho"

# Set the date as a raw serial number and reuse the existing date
# number format from the row above so no new style gets created.
$ws.Range("D6").Value = 45767
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat()

# Move the active selection to the next empty cell in the Fecha column,
# matching where Excel would leave the cursor after the edit.
$null = $ws.Range("D7").Select()
